$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.916.57'
$ws.Range('D3').Value = '1.640.92'
$ws.Range('E3').Value = '  +1.19%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.54'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.523'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.10%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.77'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.26%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.263'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.59%  '
$ws.Range('E10').Value = '  +0.79%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0876'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.47%  '
$ws.Range('D12').Value = '1.874.15'
$ws.Range('E12').Value = '  +1.21%  '
$ws.Range('D13').Value = '1.644.91'
$ws.Range('E13').Value = '  +1.17%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.09'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.12%  '
$ws.Range('E15').Value = '  +4.35%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.01'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.27%  '
$ws.Range('D17').Value = '27.912.56'
$ws.Range('E17').Value = '  +1.58%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '230.82'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.45%  '
$ws.Range('E19').Value = '  +1.28%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.60'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.99%  '
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.18'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +7.58%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('E23').Value = '  +1.52%  '
$ws.Range('E24').Value = '  -1.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.02'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.93%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.92'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.79%  '
$ws.Range('E27').Value = '  +0.80%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.72'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.20%  '
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('E30').Value = '  +1.18%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0485'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.53%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.33'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.89%  '
$ws.Range('D33').Value = '1.422.22'
$ws.Range('E33').Value = '  -2.96%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.11'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.22%  '
$ws.Range('E35').Value = '  +1.85%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.34'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.29%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.889'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.08%  '
$ws.Range('E38').Value = '  +0.83%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.924'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.66%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.557'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.96%  '
$ws.Range('E41').Value = '  +2.27%  '
$ws.Range('E42').Value = '  -0.05%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '66.96'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.46%  '
$ws.Range('E44').Value = '  +0.44%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.82'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.97%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.44'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.86%  '
$ws.Range('E47').Value = '  +0.30%  '
$ws.Range('D48').Value = '1.783.06'
$ws.Range('E48').Value = '  +1.23%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '88.84'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.86%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0104'
$ws.Range('E50').Value = '  -1.36%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.101'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.41%  '
